$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-09 12:01:59", 0.0008),
    @("2023-12-09 12:02:16", 0.0004),
    @("2023-12-09 12:02:29", 0.0004)
)

$startRow = 140
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
